$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "이노그리드" entry entirely (old row 12) - shifts subsequent rows up by one.
$ws.Rows.Item(12).Delete()

# 2. "이노스페이스" (now row 14) had its 확정공모가 (confirmed price) finalized to 43,300.
#    The source column stores this as text, so prefix with an apostrophe to force text entry
#    (otherwise Excel auto-converts a pure numeric string into a number), then clear the
#    resulting "number stored as text" formatting so the cell carries no extra style.
$ws.Range("D14").Value = "'43300"
$ws.Range("D14").Style = "Normal"

# 3. Updated subscription window / confirmed pricing data for 에이치엠씨아이비스팩7호 (now row 19).
$ws.Range("B19").Value = "2024.06.03~06.10"
$ws.Range("C19").Value = "7,500~9,000"
$ws.Range("D19").Value = "'11500"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = 23000
$ws.Range("F19").Value = "한국투자증권"

# 4. Updated subscription window / confirmed pricing data for 에스오에스랩 (now row 20).
$ws.Range("B20").Value = "2024.06.03~06.04"
$ws.Range("C20").Value = "2,000~2,000"
$ws.Range("D20").Value = "'2000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = 9500
$ws.Range("F20").Value = "미래에셋증권"

# 5. Append new entry for 미래에셋비전스팩5호 (row 21).
$ws.Range("A21").Value = "미래에셋비전스팩5호"
$ws.Range("B21").Value = "2024.06.03~06.04"
$ws.Range("C21").Value = "2,000~2,000"
$ws.Range("D21").Value = "'2000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = 9500
$ws.Range("F21").Value = "미래에셋증권"
